# Actualiza base de datos EC y agrega parte 1 de nuevos estado de cuenta
#
# The worker detail table (rows 16-22, cols B:G) is re-grouped: instead of
# being ordered by worker, it is now ordered by period (1605 first, then
# 1612), and a new set of workers (EDUARD / SAMIR / HUMBERTO) gets its
# outstanding-balance rows written out, while YINA keeps only her single
# 1612 row.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$rows = @(
  @("CC", "1047459226", "EDUARD ENRIQUE PARRA RODRIGUEZ", "1605", 27600, 644350),
  @("CC", "1047468957", "SAMIR ENRIQUE PARRA RODRIGUEZ",  "1605", 27600, 644350),
  @("CC", "73160872",   "HUMBERTO PARRA POLO",             "1605", 27600, 644350),
  @("CC", "45563280",   "YINA PATRICIA RODRIGUEZ IRIARTE", "1612", 25774, 644350),
  @("CC", "1047459226", "EDUARD ENRIQUE PARRA RODRIGUEZ", "1612", 25774, 644350),
  @("CC", "1047468957", "SAMIR ENRIQUE PARRA RODRIGUEZ",  "1612", 25774, 644350),
  @("CC", "73160872",   "HUMBERTO PARRA POLO",             "1612", 25774, 644350)
)

$startRow = 16
for ($i = 0; $i -lt $rows.Length; $i++) {
    $r = $startRow + $i
    $data = $rows[$i]
    $ws.Range("B$r").Value = $data[0]
    $ws.Range("C$r").Value = $data[1]
    $ws.Range("D$r").Value = $data[2]
    $ws.Range("E$r").Value = $data[3]
    $ws.Range("F$r").Value = $data[4]
    $ws.Range("G$r").Value = $data[5]
}

# Column C (N Doc Trabajador) now best-fits a narrower set of ID numbers,
# so the column shrinks. The headless engine quantizes ColumnWidth to
# 1/6-character steps, so this lands as close as possible to the 10.8163
# (char-width) target that real Excel's font metrics produced.
$ws.Columns.Item(3).ColumnWidth = 9.8

# The logo picture is anchored with "move and size with cells" inside
# columns B:C, so once column C narrows the picture is nudged left to
# keep its same visual spot relative to the table. Re-home it at the
# exact EMU position/size captured by the author's edit.
$shp = $ws.Shapes.Item(1)
$shp.Left = 750450 / 12700
$shp.Width = 975600 / 12700
$shp.Height = 612000 / 12700
